# New CodeSystem version: publish status moves from draft to active, with a
# refreshed publish date, and the previously-blank Experimental / Case
# Sensitive metadata rows get their explicit (text) values filled in.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# Status: draft -> active
$ws.Range("B6").Value = "active"

# Experimental: (empty) -> "false", written as literal text (not a boolean).
# A leading apostrophe forces text entry; then we copy the plain-text
# formatting from a sibling "value" cell (B15) back over it so the cell
# keeps its original (non quote-prefixed) style.
$ws.Range("B7").Value = "'false"
$ws.Range("B15").Copy()
$ws.Range("B7").PasteSpecial(-4122)

# Date: updated publish date/time
$ws.Range("B8").Value = "2024-12-16T14:50:05-03:00"

# Case Sensitive: (empty) -> "true", same text-not-boolean treatment.
$ws.Range("B17").Value = "'true"
$ws.Range("B15").Copy()
$ws.Range("B17").PasteSpecial(-4122)

$excel.CutCopyMode = $false
